# Updates odds/stats values on Sheet1 for the 2024-12-10 FlashScore weekly
# games export. Targets specific cells across rows 2-14 (header row 1 is
# untouched) to refresh odds figures pulled from the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.75
$ws.Range("H2").Value = 4.2
$ws.Range("I2").Value = 1.65
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 17
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1.57
$ws.Range("R2").Value = 2.35
$ws.Range("AC2").Value = 17
$ws.Range("AS2").Value = 126
$ws.Range("BA2").Value = 23

# Row 3
$ws.Range("H3").Value = 6.25
$ws.Range("J3").Value = 1.62
$ws.Range("K3").Value = 2.88
$ws.Range("M3").Value = 1.03
$ws.Range("N3").Value = 17
$ws.Range("O3").Value = 1.13
$ws.Range("P3").Value = 6
$ws.Range("Q3").Value = 1.44
$ws.Range("R3").Value = 2.7
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 1.73
$ws.Range("W3").Value = 9
$ws.Range("X3").Value = 7
$ws.Range("Y3").Value = 9.5
$ws.Range("AB3").Value = 26
$ws.Range("AC3").Value = 17
$ws.Range("AD3").Value = 12
$ws.Range("AG3").Value = 301
$ws.Range("AH3").Value = 29
$ws.Range("AJ3").Value = 34
$ws.Range("AK3").Value = 151
$ws.Range("AM3").Value = 67
$ws.Range("AR3").Value = 34
$ws.Range("AS3").Value = 101
$ws.Range("AV3").Value = 51
$ws.Range("AY3").Value = 41
$ws.Range("BA3").Value = 201
$ws.Range("BC3").Value = 301

# Row 4
$ws.Range("Q4").Value = 1.75
$ws.Range("R4").Value = 2.05

# Row 5
$ws.Range("G5").Value = 1.67
$ws.Range("I5").Value = 5.5
$ws.Range("J5").Value = 2.38
$ws.Range("L5").Value = 6.5
$ws.Range("Q5").Value = 2.63
$ws.Range("R5").Value = 1.5
$ws.Range("S5").Value = 1.57
$ws.Range("T5").Value = 2.25
$ws.Range("U5").Value = 2.5
$ws.Range("V5").Value = 1.5
$ws.Range("Y5").Value = 9.5
$ws.Range("AM5").Value = 67
$ws.Range("AO5").Value = 9
$ws.Range("AS5").Value = 301
$ws.Range("AT5").Value = 2.25
$ws.Range("BD5").Value = 151

# Row 6
$ws.Range("AW6").Value = 151

# Row 7
$ws.Range("G7").Value = 4.5
$ws.Range("H7").Value = 3.9
$ws.Range("I7").Value = 1.7
$ws.Range("J7").Value = 4.75
$ws.Range("L7").Value = 2.3
$ws.Range("O7").Value = 1.25
$ws.Range("P7").Value = 4
$ws.Range("U7").Value = 1.73
$ws.Range("V7").Value = 2
$ws.Range("Z7").Value = 51
$ws.Range("AD7").Value = 7.5
$ws.Range("AE7").Value = 15
$ws.Range("AG7").Value = 201
$ws.Range("AI7").Value = 8.5
$ws.Range("AK7").Value = 13
$ws.Range("AN7").Value = 6.5
$ws.Range("AQ7").Value = 81
$ws.Range("AR7").Value = 101
$ws.Range("AS7").Value = 201
$ws.Range("AU7").Value = 8
$ws.Range("AX7").Value = 3.75
$ws.Range("BB7").Value = 41

# Row 8
$ws.Range("AA8").Value = 26

# Row 9
$ws.Range("BD9").Value = 126

# Row 10
$ws.Range("G10").Value = 1.9
$ws.Range("I10").Value = 4.1
$ws.Range("J10").Value = 2.6
$ws.Range("K10").Value = 2.1
$ws.Range("L10").Value = 4.75
$ws.Range("N10").Value = 8.5
$ws.Range("Z10").Value = 15
$ws.Range("AI10").Value = 21
$ws.Range("AJ10").Value = 15

# Row 11
$ws.Range("G11").Value = 1.85
$ws.Range("H11").Value = 3.5
$ws.Range("Q11").Value = 1.95
$ws.Range("R11").Value = 1.9
$ws.Range("S11").Value = 1.4
$ws.Range("T11").Value = 2.75
$ws.Range("U11").Value = 1.8
$ws.Range("V11").Value = 1.91
$ws.Range("W11").Value = 7.5
$ws.Range("AB11").Value = 26
$ws.Range("AC11").Value = 10
$ws.Range("AD11").Value = 6.5
$ws.Range("AG11").Value = 251
$ws.Range("AM11").Value = 41
$ws.Range("AN11").Value = 3.75
$ws.Range("AP11").Value = 21
$ws.Range("AS11").Value = 151
$ws.Range("AT11").Value = 2.75
$ws.Range("BA11").Value = 81
$ws.Range("BB11").Value = 101

# Row 12
$ws.Range("H12").Value = 3.5
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 2.15
$ws.Range("K12").Value = 2.2
$ws.Range("L12").Value = 5
$ws.Range("M12").Value = 1.09
$ws.Range("N12").Value = 6.81
$ws.Range("O12").Value = 1.32
$ws.Range("P12").Value = 2.82
$ws.Range("Q12").Value = 1.98
$ws.Range("R12").Value = 1.65
$ws.Range("S12").Value = 1.39
$ws.Range("T12").Value = 2.55
$ws.Range("U12").Value = 1.91
$ws.Range("V12").Value = 1.7
$ws.Range("W12").Value = 5.9
$ws.Range("X12").Value = 7.1
$ws.Range("Y12").Value = 8.25
$ws.Range("AA12").Value = 14.5
$ws.Range("AB12").Value = 30
$ws.Range("AC12").Value = 8.75
$ws.Range("AD12").Value = 6.9
$ws.Range("AE12").Value = 17.5
$ws.Range("AG12").Value = 800
$ws.Range("AM12").Value = 60
$ws.Range("AN12").Value = 3.45
$ws.Range("AO12").Value = 7.7
$ws.Range("AP12").Value = 16.5
$ws.Range("AQ12").Value = 25
$ws.Range("AR12").Value = 55
$ws.Range("AS12").Value = 200
$ws.Range("AT12").Value = 2.65
$ws.Range("AU12").Value = 7.3
$ws.Range("AV12").Value = 65
$ws.Range("AX12").Value = 6.7
$ws.Range("AY12").Value = 28
$ws.Range("AZ12").Value = 30
$ws.Range("BA12").Value = 175
$ws.Range("BB12").Value = 175
$ws.Range("BC12").Value = 400

# Row 13
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

# Row 14
$ws.Range("G14").Value = 2.5
$ws.Range("H14").Value = 3.2
$ws.Range("I14").Value = 2.6
$ws.Range("J14").Value = 3
$ws.Range("L14").Value = 3.2
$ws.Range("P14").Value = 3.8
$ws.Range("T14").Value = 2.92
$ws.Range("U14").Value = 1.53
$ws.Range("W14").Value = 11
$ws.Range("X14").Value = 15.5
$ws.Range("Y14").Value = 9.25
$ws.Range("Z14").Value = 30
$ws.Range("AA14").Value = 18.5
$ws.Range("AD14").Value = 6.5
$ws.Range("AH14").Value = 10
$ws.Range("AI14").Value = 14.5
$ws.Range("AJ14").Value = 9.5
$ws.Range("AK14").Value = 32
$ws.Range("AL14").Value = 21
$ws.Range("AN14").Value = 4.7
$ws.Range("AO14").Value = 13
$ws.Range("AP14").Value = 17.5
$ws.Range("AQ14").Value = 50
$ws.Range("AR14").Value = 70
$ws.Range("AT14").Value = 2.92
$ws.Range("AX14").Value = 4.8
$ws.Range("AY14").Value = 14.5
$ws.Range("AZ14").Value = 19.5
$ws.Range("BA14").Value = 60
